{"js": "// Replace the body of the second paragraph (the introductory paragraph about\n// global warming / carbon footprint) with the new, rewritten text while\n// keeping the existing paragraph formatting (font, size, indentation, etc.)\n// untouched: Word.InsertLocation.replace swaps only the paragraph's text\n// content/runs, not its paragraph-level properties.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[1];\n\nconst newText =\n  \"Today, just about every human activity contributes to the increasing issue of climate change and global warming.\" +\n  \" All human activities to some extent, have direct and indirect impacts on the environment, through the emission of greenhouse gases.\" +\n  \" The term \\u201Ccarbon footprint\\u201D, is a metaphor used to describe the environmental impact of an entity or activity by estimating its total carbon emissions caused.\" +\n  \" \\u201CCarbon\\u201D, does not only refer to carbon dioxide, and includes other greenhouses gases such as methane and nitrous oxide.\" +\n  \" A carbon footprint does not only consist of direct emissions caused, but also indirect emissions. Direct emissions would refer to the combined emissions of activities that constitute the footprint.\" +\n  \" On the other hand, indirect emissions are the emissions caused by activities that lead up to the activities that constitute the footprint.\" +\n  \" For example, the direct emissions of creating a plastic box would include the emissions of all activities happening at the factory where the construction took place. Indirect emissions would include the emissions of the refinery process required to create the oil that supplied energy to the machines that created the box, and so on.\" +\n  \" In the end, being completely accurate in measuring a carbon footprint is near impossible due to the vast number of factors involved.\";\n\ntargetParagraph.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the body of the second paragraph (the introductory paragraph about\n# global warming / carbon footprint) with the new, rewritten text, keeping\n# the paragraph's own formatting (font, size, indentation, spacing) intact.\n# We trim the trailing paragraph mark from the Range before assigning new\n# text so the paragraph break itself is not disturbed.\n\n$d = $word.ActiveDocument\n\n$newText = \"Today, just about every human activity contributes to the increasing issue of climate change and global warming.\" + `\n  \" All human activities to some extent, have direct and indirect impacts on the environment, through the emission of greenhouse gases.\" + `\n  \" The term \" + [char]0x201C + \"carbon footprint\" + [char]0x201D + \", is a metaphor used to describe the environmental impact of an entity or activity by estimating its total carbon emissions caused.\" + `\n  \" \" + [char]0x201C + \"Carbon\" + [char]0x201D + \", does not only refer to carbon dioxide, and includes other greenhouses gases such as methane and nitrous oxide.\" + `\n  \" A carbon footprint does not only consist of direct emissions caused, but also indirect emissions. Direct emissions would refer to the combined emissions of activities that constitute the footprint.\" + `\n  \" On the other hand, indirect emissions are the emissions caused by activities that lead up to the activities that constitute the footprint.\" + `\n  \" For example, the direct emissions of creating a plastic box would include the emissions of all activities happening at the factory where the construction took place. Indirect emissions would include the emissions of the refinery process required to create the oil that supplied energy to the machines that created the box, and so on.\" + `\n  \" In the end, being completely accurate in measuring a carbon footprint is near impossible due to the vast number of factors involved.\"\n\n$p = $d.Paragraphs.Item(2)\n$r = $p.Range\n$r.End = $r.End - 1\n$r.Text = $newText\n"}
